# Apply the tracked changes to "tracking avancement.xlsx":
#  - Mark a batch of checklist cells (columns F/G/H, a few rows also C/D/E)
#    as completed (TRUE) on the "Feuil1" sheet.
#  - Update the saved view state (frozen-pane scroll position + active
#    selection) to reflect where the user left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- Checklist cells flipped from FALSE to TRUE -----------------------
# Rows 27-31 & 40-49: only columns G/H toggled.
$rowsGH = 27, 28, 29, 30, 31, 40, 41, 42, 43, 45, 46, 47, 49
foreach ($r in $rowsGH) {
    $ws.Range("G$r").Value = $true
    $ws.Range("H$r").Value = $true
}

# Rows 33-36: columns F/G/H toggled.
$rowsFGH = 33, 34, 35, 36
foreach ($r in $rowsFGH) {
    $ws.Range("F$r").Value = $true
    $ws.Range("G$r").Value = $true
    $ws.Range("H$r").Value = $true
}

# Row 38: the whole C:H block toggled.
$ws.Range("C38:H38").Value = $true

# --- Saved view state ---------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G30").Select() | Out-Null
